# civ_pretas_resultats.xlsx -- "Some improvements on ivory coast forms."
#
# Changes applied (derived from the target diff):
#  - survey!row5/row6 : code_genere/code_genere2 -> code_id/code_id2, type string ->
#    integer, hint/constraint/constraint_message text updated (6-digit id instead
#    of the 3-4 digit generated code).
#  - survey!C10..C13  : renumber the question labels ("Description d'invalidation..."
#    -> "5a. ...", "Si Autre (expliquer)" -> "5b. ...", "7. Photos" -> "6. Photos",
#    "8. Defis operationnels observes" -> "7. ...").
#  - survey           : a new row/question "final_result" is inserted right after
#    "defis_operationnels_autre" and before "observations"; it becomes question 8
#    ("Interpretation du resultat final des deux tests").
#  - settings         : form_title/form_id bumped from V3 to V4, version bumped
#    20200310 -> 20200320.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# 1. survey: row 5 - code_genere -> code_id
# ---------------------------------------------------------------------------
$survey.Range("A5").Value = 'integer'
$survey.Range("B5").Value = 'code_id'
$survey.Range("C5").Value = '1.a. Saisir le code d''identification Ce code doit aussi être inscrit dans le test de diagnostic rapide tels que le FTS'
$survey.Range("D5").Value = 'Exemple "190005"'
$survey.Range("F5").Value = 'regex(., ''^[0-9]{6}$'')'
$survey.Range("G5").Value = 'Le format est incorrect.'

# ---------------------------------------------------------------------------
# 2. survey: row 6 - code_genere2 -> code_id2
# ---------------------------------------------------------------------------
$survey.Range("A6").Value = 'integer'
$survey.Range("B6").Value = 'code_id2'
$survey.Range("C6").Value = '1.b. Répéter le code d''identification'
$survey.Range("D6").Value = 'Exemple "190005"'
$survey.Range("F6").Value = '. = ${code_id}'
$survey.Range("G6").Value = 'Le code n''est pas le même'

# ---------------------------------------------------------------------------
# 3. survey: renumber labels on rows 10-13
# ---------------------------------------------------------------------------
$survey.Range("C10").Value = '5a. Description d''invalidation de test 1'
$survey.Range("C11").Value = '5b. Si Autre (expliquer)'
$survey.Range("C12").Value = '6. Photos'
$survey.Range("C13").Value = '7. Défis opérationnels observés'

# ---------------------------------------------------------------------------
# 4. survey: insert the new "final_result" question as row 15, pushing the
#    rest (observations/start/end/blank formatting row) down by one.
# ---------------------------------------------------------------------------
$survey.Rows.Item(15).Insert()

# the row-insert carries column formatting down from row 14 into the blank
# cells of the new row; clear it so the row starts out unstyled.
$survey.Range("A15:L15").ClearFormats()

$survey.Range("A15").Value = 'select_one resultat_list'
$survey.Range("B15").Value = 'final_result'
$survey.Range("C15").Value = '8. Interprétation du résultat final des deux tests'
$survey.Range("H15").Value = '${num_test} = ''2'''
$survey.Range("J15").Value = 'yes'

# restore the wrap-text style on C15 (same style as the other question-label
# cells, e.g. C11) by copying its format over.
$survey.Range("C11").Copy() | Out-Null
$survey.Range("C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# keep the sheet's active cell in sync with the new layout
$survey.Range("B15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. settings: bump form_title / form_id / version from V3 to V4
# ---------------------------------------------------------------------------
$settings.Range("A2").Value = '3. Côte d''Ivoire - Pré TAS FL Résultats V4'
$settings.Range("B2").Value = 'ci_pretas_lf_resultats_v4'
$settings.Range("C2").Value = 20200320
